$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (Sizes row), pushing it down to row 3.
$ws.Rows("2:2").Insert()

# New row 2: "#My colors" sub-header + translated color values
$ws.Range("A2").Value = "#My colors"
$ws.Range("B2").Value = "Red!"
$ws.Range("C2").Value = "Green!"
$ws.Range("D2").Value = "Blue!"

# New row 4 (appended after the now-shifted Sizes row 3): "#My sizes" sub-header + translated size values
$ws.Range("A4").Value = "#My sizes"
$ws.Range("B4").Value = "small~"
$ws.Range("C4").Value = "medium~"
$ws.Range("D4").Value = "large~"

# Row 4 needs the same style as the rest of the table (row 2 inherited it automatically via Insert)
$ws.Range("A1:D1").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths (auto-fit-like values captured by Excel after editing; the
# inputs below are chosen so the engine's internal pixel rounding lands on
# the closest achievable value to the authored widths 11 / 7.25 / 9.83203125 / 7.08203125)
$ws.Columns("A").ColumnWidth = 10.285714285714286
$ws.Columns("B").ColumnWidth = 6.571428571428571
$ws.Columns("C").ColumnWidth = 9.142857142857142
$ws.Columns("D").ColumnWidth = 6.428571428571429

# Update selection to match final cursor position
$null = $ws.Range("E6").Select()
